$d = $word.ActiveDocument

$replacements = @(
    @("2023-07-30 Sunday", "2023-07-31 Monday"),
    @("61÷2=", "18÷9="),
    @("41÷4=", "91÷5="),
    @("86÷2=", "14÷5="),
    @("44÷5=", "84÷8="),
    @("15÷5=", "53÷9="),
    @("98÷3=", "65÷9="),
    @("94÷7=", "23÷7="),
    @("27÷6=", "30÷6="),
    @("77÷2=", "83÷6="),
    @("70÷3=", "44÷4="),
    @("91÷4=", "13÷3="),
    @("59÷4=", "42÷6="),
    @("65÷2=", "25÷8="),
    @("74÷9=", "76÷5="),
    @("29÷5=", "94÷9="),
    @("89÷3=", "87÷4="),
    @("12÷6=", "31÷8="),
    @("10÷7=", "97÷7="),
    @("47÷9=", "84÷2="),
    @("59÷5=", "23÷5="),
    @("11÷3=", "96÷3="),
    @("60÷8=", "13÷9="),
    @("94÷8=", "64÷6="),
    @("58÷9=", "96÷7="),
    @("65÷6=", "66÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
